# Daily attendance processing - 2025-12-04 07:28:57
# Applies the diff to "Session Analysis Results": reordered Recorded-By
# email lists, refreshed attendance counters / class+group statistics,
# and four sessions (rows 46, 107, 142, 164) that flipped from
# Not Recorded/Pending to Recorded with their own Recorded-By + counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write $text into $cellRef as a literal text value without
# letting Excel's "looks like a percentage" auto-conversion turn strings
# such as "20.5%" into a formatted number (which would also mint a new
# cell style). We stage the text in an unused scratch cell that's been
# forced to Text format, then copy/paste-special *values only* onto the
# destination so the destination keeps its original style untouched.
function Set-TextValue($cellRef, $text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $scratch.Clear()
}

# Helper: re-color/re-style a whole "session" row (A..I) by copying the
# formats from a known-good row that already has the target look
# ("Recorded" = green, style index 2), without touching its values.
function Copy-RowFormat($srcRowRange, $destRowRange) {
    $srcRowRange.Copy()
    $destRowRange.PasteSpecial(-4122)
}

# ---------------------------------------------------------------------
# Row 2 / Row 24: ANATOMY session 1 "Recorded By" list reordered
# ---------------------------------------------------------------------
$ws.Range("G2").Value = "rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
$ws.Range("G24").Value = "rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
$ws.Range("H24").Value = "153/217"

# ---------------------------------------------------------------------
# Class Statistics block (K6:L10)
# ---------------------------------------------------------------------
$ws.Range("L6").Value = 36
$ws.Range("L7").Value = 11
$ws.Range("L8").Value = 129
Set-TextValue "L9" "20.5%"
Set-TextValue "L10" "29.6%"

# ---------------------------------------------------------------------
# Group Statistics block (K14:S22)
# ---------------------------------------------------------------------
Set-TextValue "S16" "37.2%"

$ws.Range("O17").Value = 5
$ws.Range("P17").Value = 1
Set-TextValue "R17" "22.7%"
Set-TextValue "S17" "34.6%"

$ws.Range("O19").Value = 4
$ws.Range("P19").Value = 3
Set-TextValue "R19" "18.2%"
Set-TextValue "S19" "26.6%"

$ws.Range("O21").Value = 5
$ws.Range("P21").Value = 0
Set-TextValue "R21" "22.7%"
Set-TextValue "S21" "27.1%"

$ws.Range("O22").Value = 5
$ws.Range("Q22").Value = 16
Set-TextValue "R22" "22.7%"
Set-TextValue "S22" "13.5%"

# ---------------------------------------------------------------------
# PHYSIOLOGY sessions 1 & 2 "Recorded By" lists reordered (rows 18/19,
# 40/41, 150/172) — rows 19/150/172 also gained Monica.Eshak
# ---------------------------------------------------------------------
$physio1 = "aya.hanafy@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$physio2 = "Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"

$ws.Range("G18").Value = $physio1
$ws.Range("G19").Value = $physio2
$ws.Range("G40").Value = $physio1
$ws.Range("G41").Value = $physio2
$ws.Range("G150").Value = $physio2
$ws.Range("H150").Value = "95/224"
$ws.Range("G172").Value = $physio2

# ---------------------------------------------------------------------
# Other "Recorded By" email-list reorders (no count/status changes)
# ---------------------------------------------------------------------
$pharm1 = "merna.said@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, maimustafa@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg"
$ws.Range("G54").Value = $pharm1
$ws.Range("G76").Value = $pharm1
$ws.Range("G98").Value = $pharm1
$ws.Range("G120").Value = $pharm1

$pharm2 = "afaf.abdallah@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg"
$ws.Range("G58").Value = $pharm2
$ws.Range("G80").Value = $pharm2

$physioB = "yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"
$ws.Range("G62").Value = $physioB
$ws.Range("G84").Value = $physioB

$histo = "mariam.noureldin@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg"
$ws.Range("G96").Value = $histo
$ws.Range("G118").Value = $histo

$microB = "neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
$ws.Range("G106").Value = $microB
$ws.Range("G128").Value = $microB

$ws.Range("G134").Value = "majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G156").Value = "alshimaa.atef@med.asu.edu.egm, majorelle.magdy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"

# ---------------------------------------------------------------------
# Rows 46 / 107 / 142 / 164: sessions moved from Not Recorded/Pending to
# Recorded. Re-color the row to match the "Recorded" (green) look used
# elsewhere (row 2) and fill in the Recorded By / Students / Status.
# ---------------------------------------------------------------------
Copy-RowFormat $ws.Range("A2:I2") $ws.Range("A46:I46")
$ws.Range("G46").Value = "shaimaa.ahmed@med.asu.edu.eg"
$ws.Range("H46").Value = "22/220"
$ws.Range("I46").Value = "Recorded"

Copy-RowFormat $ws.Range("A2:I2") $ws.Range("A107:I107")
$ws.Range("G107").Value = "Monica.Eshak@med.asu.edu.eg"
$ws.Range("H107").Value = "9/154"
$ws.Range("I107").Value = "Recorded"

Copy-RowFormat $ws.Range("A2:I2") $ws.Range("A142:I142")
$ws.Range("G142").Value = "basma.hamed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"
$ws.Range("H142").Value = "62/224"
$ws.Range("I142").Value = "Recorded"

Copy-RowFormat $ws.Range("A2:I2") $ws.Range("A164:I164")
$ws.Range("G164").Value = "basma.hamed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"
$ws.Range("H164").Value = "1/226"
$ws.Range("I164").Value = "Recorded"
